$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1133
$ws.Range("I40").Value = 1114.4546
$ws.Range("J40").Value = 1162.1428
$ws.Range("K40").Value = 1114.4546
$ws.Range("L40").Value = 1162.1428
$ws.Range("M40").Value = -939.4546
$ws.Range("N40").Value = -1512.1428
$ws.Range("H64").Value = 3998.074
$ws.Range("J64").Value = 4073
$ws.Range("L64").Value = 4073
$ws.Range("N64").Value = -4569
$ws.Range("H67").Value = 3998.074
$ws.Range("J67").Value = 4073
$ws.Range("L67").Value = 4073
$ws.Range("N67").Value = -5789
$ws.Range("H69").Value = 1547.4073
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 1518.4615
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 4555.3845
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -6303.3845
$ws.Range("H72").Value = 1547.4073
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 1518.4615
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 13666.1535
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -22402.1535
$ws.Range("H74").Value = 4285
$ws.Range("I74").Value = 3733.3333
$ws.Range("K74").Value = 3733.3333
$ws.Range("M74").Value = -2797.3333
$ws.Range("H77").Value = 4285
$ws.Range("I77").Value = 3733.3333
$ws.Range("K77").Value = 18666.6665
$ws.Range("M77").Value = -13986.6665
$ws.Range("H97").Value = 1555
$ws.Range("J97").Value = 1555
$ws.Range("L97").Value = 4665
$ws.Range("N97").Value = -5657
$ws.Range("H135").Value = 27789640
$ws.Range("I135").Value = 1079.1818
$ws.Range("J135").Value = 71457380
$ws.Range("K135").Value = 9712.636200000001
$ws.Range("L135").Value = 643116420
$ws.Range("M135").Value = -7177.636200000001
$ws.Range("N135").Value = -643121490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2178.1555
$ws.Range("I32").Value = 1902.579
$ws.Range("K32").Value = 1902.579
$ws.Range("M32").Value = -1615.579
$ws.Range("H61").Value = 1002811.06
$ws.Range("I61").Value = 1801499.9
$ws.Range("J61").Value = 4450
$ws.Range("K61").Value = 1801499.9
$ws.Range("L61").Value = 4450
$ws.Range("M61").Value = -1801287.9
$ws.Range("N61").Value = -4874
$ws.Range("H74").Value = 3435.6667
$ws.Range("I74").Value = 2927.625
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 2927.625
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = -2053.625
$ws.Range("N74").Value = -9248
$ws.Range("H77").Value = 3435.6667
$ws.Range("I77").Value = 2927.625
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 14638.125
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = -10270.125
$ws.Range("N77").Value = -46236
$ws.Range("H132").Value = 23594.209
$ws.Range("I132").Value = 2558.3333
$ws.Range("J132").Value = 58654
$ws.Range("K132").Value = 7674.999899999999
$ws.Range("L132").Value = 175962
$ws.Range("M132").Value = -5144.999899999999
$ws.Range("N132").Value = -181022
$ws.Range("H136").Value = 1002811.06
$ws.Range("I136").Value = 1801499.9
$ws.Range("J136").Value = 4450
$ws.Range("K136").Value = 5404499.699999999
$ws.Range("L136").Value = 13350
$ws.Range("M136").Value = -5401949.699999999
$ws.Range("N136").Value = -18450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2052.1667
$ws.Range("I86").Value = 1826.0769
$ws.Range("J86").Value = 2640
$ws.Range("K86").Value = 1826.0769
$ws.Range("L86").Value = 2640
$ws.Range("M86").Value = -703.0769
$ws.Range("N86").Value = -4886
$ws.Range("H89").Value = 2052.1667
$ws.Range("I89").Value = 1826.0769
$ws.Range("J89").Value = 2640
$ws.Range("K89").Value = 9130.3845
$ws.Range("L89").Value = 13200
$ws.Range("M89").Value = -3514.3845
$ws.Range("N89").Value = -24432
$ws.Range("H105").Value = 2059.9333
$ws.Range("I105").Value = 2024.875
$ws.Range("K105").Value = 2024.875
$ws.Range("M105").Value = -277.875
$ws.Range("H134").Value = 5125.0586
$ws.Range("I134").Value = 5528.467
$ws.Range("K134").Value = 16585.401
$ws.Range("M134").Value = -14050.401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 34717.535
$ws.Range("I58").Value = 1844.3334
$ws.Range("J58").Value = 84027.336
$ws.Range("K58").Value = 1844.3334
$ws.Range("L58").Value = 84027.336
$ws.Range("M58").Value = -1641.3334
$ws.Range("N58").Value = -84433.336
$ws.Range("H105").Value = 41667336
$ws.Range("I105").Value = 125000000
$ws.Range("J105").Value = 1005.5
$ws.Range("K105").Value = 125000000
$ws.Range("L105").Value = 1005.5
$ws.Range("M105").Value = -124998253
$ws.Range("N105").Value = -4499.5
$ws.Range("H132").Value = 29256.37
$ws.Range("I132").Value = 42904.082
$ws.Range("K132").Value = 128712.246
$ws.Range("M132").Value = -126182.246
$ws.Range("H134").Value = 1407.7778
$ws.Range("I134").Value = 1164
$ws.Range("K134").Value = 3492
$ws.Range("M134").Value = -957
$ws.Range("H136").Value = 34717.535
$ws.Range("I136").Value = 1844.3334
$ws.Range("J136").Value = 84027.336
$ws.Range("K136").Value = 5533.0002
$ws.Range("L136").Value = 252082.008
$ws.Range("M136").Value = -2983.0002
$ws.Range("N136").Value = -257182.008

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 1000
$ws.Range("K8").Value = 3000
$ws.Range("M8").Value = -2861
$ws.Range("H107").Value = 12760.875
$ws.Range("I107").Value = 100000
$ws.Range("J107").Value = 298.14285
$ws.Range("K107").Value = 300000
$ws.Range("L107").Value = 894.4285500000001
$ws.Range("M107").Value = -298080
$ws.Range("N107").Value = -4734.428550000001
$ws.Range("H122").Value = 697.25
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 697.25
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6275.25
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11175.25
$ws.Range("H131").Value = 182610.38
$ws.Range("I131").Value = 807.5
$ws.Range("J131").Value = 196869.44
$ws.Range("K131").Value = 2422.5
$ws.Range("L131").Value = 590608.3200000001
$ws.Range("M131").Value = 2617.5
$ws.Range("N131").Value = -600688.3200000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 31689.945
$ws.Range("I132").Value = 5525.5713
$ws.Range("J132").Value = 48340
$ws.Range("K132").Value = 16576.7139
$ws.Range("L132").Value = 145020
$ws.Range("M132").Value = -14046.7139
$ws.Range("N132").Value = -150080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5680.2
$ws.Range("I68").Value = 2800.3333
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2800.3333
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -2051.3333
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 5680.2
$ws.Range("I71").Value = 2800.3333
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 14001.6665
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -10257.6665
$ws.Range("N71").Value = -57488
$ws.Range("H132").Value = 4665.8887
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5997.6665
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 17992.9995
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -23052.9995
$ws.Range("H136").Value = 3553.4614
$ws.Range("I136").Value = 2466.111
$ws.Range("K136").Value = 7398.333
$ws.Range("M136").Value = -4848.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1748.7142
$ws.Range("I81").Value = 1748.7142
$ws.Range("K81").Value = 3497.4284
$ws.Range("M81").Value = -2436.4284
$ws.Range("H84").Value = 1748.7142
$ws.Range("I84").Value = 1748.7142
$ws.Range("K84").Value = 17487.142
$ws.Range("M84").Value = -12183.142
$ws.Range("H132").Value = 3832.889
$ws.Range("I132").Value = 3200
$ws.Range("K132").Value = 9600
$ws.Range("M132").Value = -7070
$ws.Range("H136").Value = 1408.7778
$ws.Range("I136").Value = 1187.1538
$ws.Range("J136").Value = 1614.5714
$ws.Range("K136").Value = 3561.4614
$ws.Range("L136").Value = 4843.7142
$ws.Range("M136").Value = -1011.4614
$ws.Range("N136").Value = -9943.7142

Write-Host "Applied $(214) cell updates across 8 sheets"
